$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the latest cryptos-list scrape.
# Price values in column D are plain numeric-looking text (e.g. "0.999", "9.08")
# that must stay text, matching the source data. Assigning through .Value lets
# Excel auto-detect them as numbers, so we force text entry with a leading
# apostrophe and then strip the resulting quote-prefix/number-format styling
# with ClearFormats() so the cell ends up identical to the original (no direct
# formatting), just holding the new text value.
$updates = @(
    @{ Cell = "D2"; Value = '79.886.26'; ForceText = $True }
    @{ Cell = "E2"; Value = '  +4.54%  '; ForceText = $False }
    @{ Cell = "D3"; Value = '3.204.62'; ForceText = $True }
    @{ Cell = "E3"; Value = '  +5.05%  '; ForceText = $False }
    @{ Cell = "E4"; Value = '  -0.06%  '; ForceText = $False }
    @{ Cell = "D5"; Value = '206.07'; ForceText = $True }
    @{ Cell = "E5"; Value = '  +2.25%  '; ForceText = $False }
    @{ Cell = "D6"; Value = '636.37'; ForceText = $True }
    @{ Cell = "E6"; Value = '  +1.81%  '; ForceText = $False }
    @{ Cell = "D7"; Value = '0.999'; ForceText = $True }
    @{ Cell = "E7"; Value = '  +0.00%  '; ForceText = $False }
    @{ Cell = "D8"; Value = '0.239'; ForceText = $True }
    @{ Cell = "E8"; Value = '  +14.76%  '; ForceText = $False }
    @{ Cell = "D9"; Value = '0.584'; ForceText = $True }
    @{ Cell = "E9"; Value = '  +5.71%  '; ForceText = $False }
    @{ Cell = "D10"; Value = '3.201.47'; ForceText = $True }
    @{ Cell = "E10"; Value = '  +5.08%  '; ForceText = $False }
    @{ Cell = "D11"; Value = '0.582'; ForceText = $True }
    @{ Cell = "E11"; Value = '  +32.77%  '; ForceText = $False }
    @{ Cell = "E12"; Value = '  +2.93%  '; ForceText = $False }
    @{ Cell = "D13"; Value = '5.53'; ForceText = $True }
    @{ Cell = "E13"; Value = '  +7.13%  '; ForceText = $False }
    @{ Cell = "E14"; Value = '  +19.20%  '; ForceText = $False }
    @{ Cell = "D15"; Value = '3.785.87'; ForceText = $True }
    @{ Cell = "E15"; Value = '  +4.83%  '; ForceText = $False }
    @{ Cell = "D16"; Value = '31.93'; ForceText = $True }
    @{ Cell = "E16"; Value = '  +8.79%  '; ForceText = $False }
    @{ Cell = "D17"; Value = '79.538.47'; ForceText = $True }
    @{ Cell = "E17"; Value = '  +4.18%  '; ForceText = $False }
    @{ Cell = "D18"; Value = '3.186.96'; ForceText = $True }
    @{ Cell = "E18"; Value = '  +3.99%  '; ForceText = $False }
    @{ Cell = "D19"; Value = '14.52'; ForceText = $True }
    @{ Cell = "E19"; Value = '  +7.07%  '; ForceText = $False }
    @{ Cell = "D20"; Value = '3.04'; ForceText = $True }
    @{ Cell = "E20"; Value = '  +30.66%  '; ForceText = $False }
    @{ Cell = "D21"; Value = '9.22'; ForceText = $True }
    @{ Cell = "E21"; Value = '  +1.41%  '; ForceText = $False }
    @{ Cell = "D22"; Value = '433.46'; ForceText = $True }
    @{ Cell = "E22"; Value = '  +15.38%  '; ForceText = $False }
    @{ Cell = "D23"; Value = '5.15'; ForceText = $True }
    @{ Cell = "E23"; Value = '  +17.81%  '; ForceText = $False }
    @{ Cell = "D24"; Value = '11.29'; ForceText = $True }
    @{ Cell = "E24"; Value = '  +13.69%  '; ForceText = $False }
    @{ Cell = "D25"; Value = '3.365.57'; ForceText = $True }
    @{ Cell = "E25"; Value = '  +4.96%  '; ForceText = $False }
    @{ Cell = "D26"; Value = '77.10'; ForceText = $True }
    @{ Cell = "E26"; Value = '  +4.69%  '; ForceText = $False }
    @{ Cell = "D27"; Value = '4.75'; ForceText = $True }
    @{ Cell = "E27"; Value = '  +7.52%  '; ForceText = $False }
    @{ Cell = "E28"; Value = '  +0.15%  '; ForceText = $False }
    @{ Cell = "D29"; Value = '0.0000120'; ForceText = $True }
    @{ Cell = "E29"; Value = '  +6.66%  '; ForceText = $False }
    @{ Cell = "B30"; Value = 'InternetComputer(DFINITY)'; ForceText = $False }
    @{ Cell = "C30"; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $False }
    @{ Cell = "D30"; Value = '9.08'; ForceText = $True }
    @{ Cell = "E30"; Value = '  +9.12%  '; ForceText = $False }
    @{ Cell = "B31"; Value = 'Binance-PegBSC-USD'; ForceText = $False }
    @{ Cell = "C31"; Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; ForceText = $False }
    @{ Cell = "D31"; Value = '0.999'; ForceText = $True }
    @{ Cell = "E31"; Value = '  -0.15%  '; ForceText = $False }
    @{ Cell = "E32"; Value = '  +5.02%  '; ForceText = $False }
    @{ Cell = "D33"; Value = '526.35'; ForceText = $True }
    @{ Cell = "E33"; Value = '  +3.79%  '; ForceText = $False }
    @{ Cell = "E34"; Value = '  +1.81%  '; ForceText = $False }
    @{ Cell = "E35"; Value = '  +23.86%  '; ForceText = $False }
    @{ Cell = "D36"; Value = '23.18'; ForceText = $True }
    @{ Cell = "E36"; Value = '  +11.24%  '; ForceText = $False }
    @{ Cell = "E37"; Value = '  +13.04%  '; ForceText = $False }
    @{ Cell = "E38"; Value = '  -0.04%  '; ForceText = $False }
    @{ Cell = "D39"; Value = '0.409'; ForceText = $True }
    @{ Cell = "E39"; Value = '  +5.49%  '; ForceText = $False }
    @{ Cell = "D40"; Value = '165.38'; ForceText = $True }
    @{ Cell = "E40"; Value = '  +1.42%  '; ForceText = $False }
    @{ Cell = "D41"; Value = '20.05'; ForceText = $True }
    @{ Cell = "E41"; Value = '  +0.10%  '; ForceText = $False }
    @{ Cell = "D42"; Value = '192.83'; ForceText = $True }
    @{ Cell = "E42"; Value = '  +0.59%  '; ForceText = $False }
    @{ Cell = "E43"; Value = '  +0.07%  '; ForceText = $False }
    @{ Cell = "D44"; Value = '5.55'; ForceText = $True }
    @{ Cell = "E44"; Value = '  +6.81%  '; ForceText = $False }
    @{ Cell = "E45"; Value = '  +4.62%  '; ForceText = $False }
    @{ Cell = "D46"; Value = '1.81'; ForceText = $True }
    @{ Cell = "E46"; Value = '  +8.16%  '; ForceText = $False }
    @{ Cell = "D47"; Value = '1.33'; ForceText = $True }
    @{ Cell = "E47"; Value = '  +3.22%  '; ForceText = $False }
    @{ Cell = "D48"; Value = '43.25'; ForceText = $True }
    @{ Cell = "E48"; Value = '  +2.42%  '; ForceText = $False }
    @{ Cell = "D49"; Value = '25.92'; ForceText = $True }
    @{ Cell = "E49"; Value = '  +15.29%  '; ForceText = $False }
    @{ Cell = "D50"; Value = '0.644'; ForceText = $True }
    @{ Cell = "E50"; Value = '  +5.26%  '; ForceText = $False }
    @{ Cell = "D51"; Value = '2.54'; ForceText = $True }
    @{ Cell = "E51"; Value = '  +2.38%  '; ForceText = $False }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $range.Value = "'" + $u.Value
        $range.ClearFormats()
    } else {
        $range.Value = $u.Value
    }
}
